$d = $word.ActiveDocument

# Remove all content from the document (text, paragraphs, proofErr markers,
# the old bookmark, etc.) leaving a single, empty trailing paragraph, which
# is the minimum any Word document must contain.
while ($d.Content.End -gt 1) {
    $d.Content.Delete()
}

# Re-create the "_GoBack" bookmark (Word always keeps this bookmark around
# to mark the last edit position) on the now-empty document content.
$d.Bookmarks.Add("_GoBack", $d.Content) | Out-Null
